$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8063316666666666
$ws.Range("N2").Value = 2.418995
$ws.Range("O2").Value = 0.1277387112198808
$ws.Range("P2").Value = 0.1277387112198808
$ws.Range("Q2").Value = 32.86074995612555
$ws.Range("R2").Value = 295.74674960513
$ws.Range("S2").Value = 0.002716410486296761
$ws.Range("T2").Value = 0.002716410486296761

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.578098999999999
$ws.Range("N3").Value = 10.734297
$ws.Range("O3").Value = 0.566840884181833
$ws.Range("P3").Value = 0.5668408841818329
$ws.Range("Q3").Value = 145.8196687764086
$ws.Range("R3").Value = 1312.377018987678
$ws.Range("S3").Value = 0.01205407904267014
$ws.Range("T3").Value = 0.01205407904267014

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.215895
$ws.Range("N4").Value = 0.647685
$ws.Range("O4").Value = 0.03420199180918047
$ws.Range("P4").Value = 0.03420199180918047
$ws.Range("Q4").Value = 8.798453421909999
$ws.Range("R4").Value = 79.18608079719
$ws.Range("S4").Value = 0.000727317884417751
$ws.Range("T4").Value = 0.000727317884417751

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.712026
$ws.Range("N5").Value = 5.136078
$ws.Range("O5").Value = 0.2712184127891059
$ws.Range("P5").Value = 0.2712184127891059
$ws.Range("Q5").Value = 69.77086555084134
$ws.Range("R5").Value = 627.937789957572
$ws.Range("S5").Value = 0.005767558898483915
$ws.Range("T5").Value = 0.005767558898483915

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8063316666666666
$ws.Range("N6").Value = 2.418995
$ws.Range("O6").Value = 0.1277387112198808
$ws.Range("P6").Value = 0.1277387112198808
$ws.Range("Q6").Value = 1362.127462126711
$ws.Range("R6").Value = 12259.1471591404
$ws.Range("S6").Value = 0.112599296325678
$ws.Range("T6").Value = 0.112599296325678

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.578098999999999
$ws.Range("N7").Value = 10.734297
$ws.Range("O7").Value = 0.566840884181833
$ws.Range("P7").Value = 0.5668408841818329
$ws.Range("Q7").Value = 6044.444378894692
$ws.Range("R7").Value = 54399.99941005222
$ws.Range("S7").Value = 0.4996596887347168
$ws.Range("T7").Value = 0.4996596887347168

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.215895
$ws.Range("N8").Value = 0.647685
$ws.Range("O8").Value = 0.03420199180918047
$ws.Range("P8").Value = 0.03420199180918047
$ws.Range("Q8").Value = 364.7091148627999
$ws.Range("R8").Value = 3282.3820337652
$ws.Range("S8").Value = 0.03014841917436652
$ws.Range("T8").Value = 0.03014841917436652

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.712026
$ws.Range("N9").Value = 5.136078
$ws.Range("O9").Value = 0.2712184127891059
$ws.Range("P9").Value = 0.2712184127891059
$ws.Range("Q9").Value = 2892.107214535306
$ws.Range("R9").Value = 26028.96493081776
$ws.Range("S9").Value = 0.2390739826555225
$ws.Range("T9").Value = 0.2390739826555225

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8063316666666666
$ws.Range("N10").Value = 2.418995
$ws.Range("O10").Value = 0.1277387112198808
$ws.Range("P10").Value = 0.1277387112198808
$ws.Range("Q10").Value = 81.4116011391072
$ws.Range("R10").Value = 732.704410251965
$ws.Range("S10").Value = 0.006729832013443023
$ws.Range("T10").Value = 0.006729832013443023

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.578098999999999
$ws.Range("N11").Value = 10.734297
$ws.Range("O11").Value = 0.566840884181833
$ws.Range("P11").Value = 0.5668408841818329
$ws.Range("Q11").Value = 361.2642051234976
$ws.Range("R11").Value = 3251.377846111478
$ws.Range("S11").Value = 0.02986364816479794
$ws.Range("T11").Value = 0.02986364816479794

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.215895
$ws.Range("N12").Value = 0.647685
$ws.Range("O12").Value = 0.03420199180918047
$ws.Range("P12").Value = 0.03420199180918047
$ws.Range("Q12").Value = 21.797925536755
$ws.Range("R12").Value = 196.181329830795
$ws.Range("S12").Value = 0.001801909986431077
$ws.Range("T12").Value = 0.001801909986431077

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.712026
$ws.Range("N13").Value = 5.136078
$ws.Range("O13").Value = 0.2712184127891059
$ws.Range("P13").Value = 0.2712184127891059
$ws.Range("Q13").Value = 172.8553938951273
$ws.Range("R13").Value = 1555.698545056146
$ws.Range("S13").Value = 0.01428896800032262
$ws.Range("T13").Value = 0.01428896800032262

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.8063316666666666
$ws.Range("N14").Value = 2.418995
$ws.Range("O14").Value = 0.1277387112198808
$ws.Range("P14").Value = 0.1277387112198808
$ws.Range("Q14").Value = 68.87100291186556
$ws.Range("R14").Value = 619.83902620679
$ws.Range("S14").Value = 0.005693172394463033
$ws.Range("T14").Value = 0.005693172394463033

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.578098999999999
$ws.Range("N15").Value = 10.734297
$ws.Range("O15").Value = 0.566840884181833
$ws.Range("P15").Value = 0.5668408841818329
$ws.Range("Q15").Value = 305.6152658206527
$ws.Range("R15").Value = 2750.537392385874
$ws.Range("S15").Value = 0.02526346823964802
$ws.Range("T15").Value = 0.02526346823964801

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.215895
$ws.Range("N16").Value = 0.647685
$ws.Range("O16").Value = 0.03420199180918047
$ws.Range("P16").Value = 0.03420199180918047
$ws.Range("Q16").Value = 18.44018508553
$ws.Range("R16").Value = 165.96166576977
$ws.Range("S16").Value = 0.001524344763965114
$ws.Range("T16").Value = 0.001524344763965114

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.712026
$ws.Range("N17").Value = 5.136078
$ws.Range("O17").Value = 0.2712184127891059
$ws.Range("P17").Value = 0.2712184127891059
$ws.Range("Q17").Value = 146.2288441660974
$ws.Range("R17").Value = 1316.059597494876
$ws.Range("S17").Value = 0.0120879032347768
$ws.Range("T17").Value = 0.0120879032347768

Write-Output "Applied updates to rows 2-17"